$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was reported for this market; it belongs
# chronologically where the old row 8 used to sit, so insert a fresh
# row there and push the existing rows (old 8..59) down to (9..60).
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new record's data.
$ws.Range("A8").Value = 2
$ws.Range("B8").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 44749
$ws.Range("D8").NumberFormat = $ws.Range("D9").NumberFormat
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 100112022
$ws.Range("G8").Value = "Arveja Verde"
$ws.Range("H8").Value = "Perfection"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 470
$ws.Range("K8").Value = 28000
$ws.Range("L8").Value = 30000
$ws.Range("M8").Value = 29064
$ws.Range("N8").Value = "$/malla 25 kilos"
$ws.Range("O8").Value = "Provincia de Limarí"
$ws.Range("P8").Value = 1163
$ws.Range("Q8").Value = 25
$ws.Range("R8").Value = "Hortaliza"
